# "My local changes before rebase"
#
# Append the new 3rd-year AIML student roster (23D31A6601 .. 23D31A6647)
# to Sheet1 as rows 49-95, mirror the header's bold styling onto the
# newly-touched E1:G1 cells, and leave the view in the state the author's
# session ended in (selection over the new rows, widened column E).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$startRow = 49
$lastRow = $startRow + 46

# Column A (roll numbers) was filled down first, in its own pass, so the
# shared-string table picks up all 47 new roll numbers contiguously.
for ($i = 1; $i -le 47; $i++) {
    $row = $startRow + ($i - 1)
    $num = "{0:D2}" -f $i
    $roll = "23D31A66" + $num
    $ws.Cells.Item($row, 1).Value = $roll
}

# Then column B (branch), reusing the existing "AIML" shared string.
for ($row = $startRow; $row -le $lastRow; $row++) {
    $ws.Cells.Item($row, 2).Value = "AIML"
}

# Then column C (year) -- "3rd" is a brand-new shared string, so it lands
# at the end of the table, after all 47 roll numbers.
for ($row = $startRow; $row -le $lastRow; $row++) {
    $ws.Cells.Item($row, 3).Value = "3rd"
}

# The header row's formatting now extends out to column G (its row span
# widens even though only A:C carry real values) -- match the bold style
# already used by A1:C1 so it reuses the same cell style, not a new one.
$ws.Range("E1:G1").Font.Bold = $true

# Column E is now part of the used range; give it the width left behind
# in the author's workbook.
$ws.Columns.Item(5).ColumnWidth = 11

# Restore the selection the author left on the new block of rows.
$ws.Range("C53:C95").Select()
